# Generalization of the "Relative Energy" column: convert from kcal/mol
# (factor 627.5095) to kJ/mol (factor 2625.5), and relabel the header to
# make the new unit explicit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header for column G.
$ws.Range("G1").Value = "Relative Energy (kJ/mol)"

# Update the relative-energy formulas to use the kJ/mol conversion factor.
$ws.Range("G2").Formula = "=(D2-`$D`$3)*2625.5"
$ws.Range("G3:G7").Formula = "=(D3-`$D`$3)*2625.5"

# Restore the active selection to the recalculated range.
$ws.Range("G2:G7").Select() | Out-Null
